# daily auto push: 2026-01-15 22:35 UTC
# Insert one new data row at row 633 (pushing the existing rows 633-674
# down to 634-675) and populate it with the new day's first reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (old rows 633..674) down by one row.
$ws.Rows.Item(633).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel date serials. Force text format first so the "YYYY/MM/DD"
# string isn't auto-converted into a date serial number, then clear the
# temporary formatting back off so the new row matches the plain
# (unstyled) look of the rest of the data rows.
$ws.Range("A633").NumberFormat = "@"
$ws.Range("A633").Value = "2026/01/16"
$ws.Range("A633").ClearFormats()

$ws.Range("B633").Value = "金"
$ws.Range("C633").Value = 3
$ws.Range("D633").Value = 41
